# Pipeline-Processor instruction-encoding sheet: widen the "funct" field
# from 3 bits to 4 bits, add the two new opcodes (ldd/std) that free up,
# and drop the now-redundant standalone "reset" opcode row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: bit-range markers (displayed as h:mm). shamt shifts up one bit
#     (was bits 4-3, now 5-4) and funct grows to 4 bits (was bits 2-0, now 3-0).
$ws.Range("F2").Value = 0.21111111111111111
$ws.Range("G2").Value = 0.125

# --- New 4-bit opcodes for ldd rs rd / std rs rd
$ws.Range("C25").Value = "1011"
$ws.Range("C26").Value = "1100"

# --- funct column (G) is now 4 bits wide instead of 3
$ws.Range("G3").Value  = "0000"
$ws.Range("G4").Value  = "1001"
$ws.Range("G5").Value  = "1010"
$ws.Range("G6").Value  = "0001"
$ws.Range("G7").Value  = "0111"
$ws.Range("G8").Value  = "1000"
$ws.Range("G9").Value  = "xxxx"
$ws.Range("G10").Value = "xxxx"
$ws.Range("G12").Value = "0011"
$ws.Range("G13").Value = "0010"
$ws.Range("G14").Value = "0100"
$ws.Range("G15").Value = "0101"
$ws.Range("G16").Value = "0110"
$ws.Range("G17").Value = "1001"
$ws.Range("G18").Value = "xxxx"
$ws.Range("G20").Value = "xxxx"
$ws.Range("G21").Value = "xxxx"
$ws.Range("G24").Value = "0011"
$ws.Range("G25").Value = "0011"
$ws.Range("G26").Value = "0011"
$ws.Range("G28").Value = "xxxx"
$ws.Range("G29").Value = "xxxx"
$ws.Range("G30").Value = "xxxx"
$ws.Range("G31").Value = "xxxx"
$ws.Range("G32").Value = "xxxx"
$ws.Range("G33").Value = "xxxx"
$ws.Range("G34").Value = "xxxx"

# --- Highlight the funct codes that were reassigned to new values in red
$redCells = @("G6","G13","G14","G15","G16","G17","G25","G26")
foreach ($addr in $redCells) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Font.Size = 11
    $rng.Font.Color = 255
}

# --- Drop the standalone "reset" funct-code row (old row 35) and clear the
#     leftover opcode value that used to sit next to "reset" (old row 36)
$ws.Range("C35").Value = $null
$ws.Range("C36").Value = $null

# --- Restore the cursor/selection position recorded for this edit
$ws.Range("G37").Select()
